$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-28 23:18:34"
$ws.Range("N2").Value = "-0.1 °C 22:48 TU"
$ws.Range("O2").Value = "2.6 °C"
$ws.Range("E3").Value = "2026-02-28 23:18:37"
$ws.Range("O3").Value = "-1.5 °C"
$ws.Range("E4").Value = "2026-02-28 23:18:40"
$ws.Range("O4").Value = "10.6 °C"
$ws.Range("E5").Value = "2026-02-28 23:18:42"
$ws.Range("E6").Value = "2026-02-28 23:18:45"
$ws.Range("O6").Value = "11.6 °C"
$ws.Range("E7").Value = "2026-02-28 23:18:47"
$ws.Range("J7").Value = "1025.0 hPa"
$ws.Range("E8").Value = "2026-02-28 23:18:50"
$ws.Range("E9").Value = "2026-02-28 23:18:53"
$ws.Range("E10").Value = "2026-02-28 23:18:56"
$ws.Range("O10").Value = "10.4 °C"
$ws.Range("E11").Value = "2026-02-28 23:18:58"
$ws.Range("I11").Value = "0.2 mm"
$ws.Range("E12").Value = "2026-02-28 23:19:01"
$ws.Range("E13").Value = "2026-02-28 23:19:03"
$ws.Range("I13").Value = "0.1 mm"
$ws.Range("J13").Value = "1024.4 hPa"
$ws.Range("O13").Value = "6.3 °C"
$ws.Range("E14").Value = "2026-02-28 23:19:06"
$ws.Range("E15").Value = "2026-02-28 23:19:07"
$ws.Range("O15").Value = "10.7 °C"
$ws.Range("E16").Value = "2026-02-28 23:19:08"
$ws.Range("N16").Value = "-3.7 °C 22:59 TU"
$ws.Range("O16").Value = "-1.4 °C"
$ws.Range("E17").Value = "2026-02-28 23:19:09"
$ws.Range("E18").Value = "2026-02-28 23:19:11"
$ws.Range("N18").Value = "6.0 °C 22:38 TU"
$ws.Range("O18").Value = "11.1 °C"
$ws.Range("E19").Value = "2026-02-28 23:19:12"
$ws.Range("I19").Value = "1.6 mm"
$ws.Range("E20").Value = "2026-02-28 23:19:13"
$ws.Range("N20").Value = "-2.3 °C 22:58 TU"
$ws.Range("E21").Value = "2026-02-28 23:19:14"
$ws.Range("E22").Value = "2026-02-28 23:19:17"
$ws.Range("N22").Value = "-4.0 °C 22:59 TU"
$ws.Range("E23").Value = "2026-02-28 23:19:19"
$ws.Range("E24").Value = "2026-02-28 23:19:22"
$ws.Range("E25").Value = "2026-02-28 23:19:24"
$ws.Range("I25").Value = "2.1 mm"
$ws.Range("O25").Value = "0.9 °C"
$ws.Range("E26").Value = "2026-02-28 23:19:27"
$ws.Range("I26").Value = "1.7 mm"
$ws.Range("N26").Value = "2.9 °C 22:59 TU"
$ws.Range("E27").Value = "2026-02-28 23:19:30"
$ws.Range("N27").Value = "-1.0 °C 22:59 TU"
$ws.Range("E28").Value = "2026-02-28 23:19:32"
$ws.Range("I28").Value = "0.5 mm"
$ws.Range("J28").Value = "1024.9 hPa"
$ws.Range("O28").Value = "9.6 °C"
$ws.Range("E29").Value = "2026-02-28 23:19:35"
$ws.Range("E30").Value = "2026-02-28 23:19:37"
$ws.Range("E31").Value = "2026-02-28 23:19:40"
$ws.Range("E32").Value = "2026-02-28 23:19:43"
$ws.Range("N32").Value = "4.1 °C 22:57 TU"
$ws.Range("O32").Value = "5.6 °C"
$ws.Range("E33").Value = "2026-02-28 23:19:45"
$ws.Range("I33").Value = "1.3 mm"
$ws.Range("J33").Value = "1023.4 hPa"
$ws.Range("E34").Value = "2026-02-28 23:19:48"
$ws.Range("I34").Value = "2.8 mm"
$ws.Range("E35").Value = "2026-02-28 23:19:51"
$ws.Range("J35").Value = "1025.2 hPa"
$ws.Range("E36").Value = "2026-02-28 23:19:53"
$ws.Range("O36").Value = "12.4 °C"
$ws.Range("E37").Value = "2026-02-28 23:19:56"
$ws.Range("I37").Value = "1.8 mm"
$ws.Range("J37").Value = "1026.0 hPa"
$ws.Range("O37").Value = "7.1 °C"
$ws.Range("E38").Value = "2026-02-28 23:19:59"
$ws.Range("O38").Value = "11.3 °C"
$ws.Range("E39").Value = "2026-02-28 23:20:02"
$ws.Range("E40").Value = "2026-02-28 23:20:04"
$ws.Range("J40").Value = "1024.7 hPa"
$ws.Range("E41").Value = "2026-02-28 23:20:07"
$ws.Range("J41").Value = "1024.8 hPa"
$ws.Range("O41").Value = "13.1 °C"
$ws.Range("E42").Value = "2026-02-28 23:20:09"
$ws.Range("O42").Value = "10.8 °C"
$ws.Range("E43").Value = "2026-02-28 23:20:11"
$ws.Range("E44").Value = "2026-02-28 23:20:14"
$ws.Range("O44").Value = "-1.3 °C"
$ws.Range("E45").Value = "2026-02-28 23:20:17"
$ws.Range("N45").Value = "3.8 °C 22:55 TU"
$ws.Range("E46").Value = "2026-02-28 23:20:20"
$ws.Range("J46").Value = "1025.3 hPa"

# Percentage-formatted text values: Excel auto-converts plain "NN%" strings
# to numeric percentages, which would also change the cell style (numFmtId).
# Force Text format, set the literal value, then restore the original
# "General" style (s=3) by pasting formats from the same-row G column cell,
# which always carries the untouched style.
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "69%"
$ws.Range("G13").Copy()
$ws.Range("H13").PasteSpecial(-4122)

$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "70%"
$ws.Range("G16").Copy()
$ws.Range("H16").PasteSpecial(-4122)

$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "68%"
$ws.Range("G20").Copy()
$ws.Range("H20").PasteSpecial(-4122)

$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "72%"
$ws.Range("G22").Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "74%"
$ws.Range("G23").Copy()
$ws.Range("H23").PasteSpecial(-4122)

$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "98%"
$ws.Range("G24").Copy()
$ws.Range("H24").PasteSpecial(-4122)

$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "67%"
$ws.Range("G33").Copy()
$ws.Range("H33").PasteSpecial(-4122)

$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "72%"
$ws.Range("G34").Copy()
$ws.Range("H34").PasteSpecial(-4122)

$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "76%"
$ws.Range("G40").Copy()
$ws.Range("H40").PasteSpecial(-4122)

$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "75%"
$ws.Range("G41").Copy()
$ws.Range("H41").PasteSpecial(-4122)

$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "82%"
$ws.Range("G43").Copy()
$ws.Range("H43").PasteSpecial(-4122)

$excel.CutCopyMode = $false
